$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 63.2  # H33: was 63.6
$ws.Cells.Item(33, 9).Value = 29  # I33: was 29.5
$ws.Cells.Item(33, 11).Value = 29  # K33: was 29.5
$ws.Cells.Item(33, 13).Value = 200  # M33: was 199.5

$ws.Cells.Item(51, 8).Value = 3554.3333  # H51: was 3624.375
$ws.Cells.Item(51, 10).Value = 3554.3333  # J51: was 3624.375
$ws.Cells.Item(51, 12).Value = 3554.3333  # L51: was 3624.375
$ws.Cells.Item(51, 14).Value = -4522.3333  # N51: was -4592.375

$ws.Cells.Item(64, 8).Value = 8450  # H64: was 7519.8
$ws.Cells.Item(64, 9).Value = 0  # I64: was 3799
$ws.Cells.Item(64, 11).Value = 0  # K64: was 3799
$ws.Cells.Item(64, 13).ClearContents()  # M64: was -3551

$ws.Cells.Item(67, 8).Value = 8450  # H67: was 7519.8
$ws.Cells.Item(67, 9).Value = 0  # I67: was 3799
$ws.Cells.Item(67, 11).Value = 0  # K67: was 3799
$ws.Cells.Item(67, 13).ClearContents()  # M67: was -2941

$ws.Cells.Item(70, 8).Value = 2987.92  # H70: was 3095.75
$ws.Cells.Item(70, 10).Value = 3756.8667  # J70: was 3996.6428
$ws.Cells.Item(70, 12).Value = 11270.6001  # L70: was 11989.9284
$ws.Cells.Item(70, 14).Value = -11810.6001  # N70: was -12529.9284

$ws.Cells.Item(73, 8).Value = 2987.92  # H73: was 3095.75
$ws.Cells.Item(73, 10).Value = 3756.8667  # J73: was 3996.6428
$ws.Cells.Item(73, 12).Value = 11270.6001  # L73: was 11989.9284
$ws.Cells.Item(73, 14).Value = -13142.6001  # N73: was -13861.9284

$ws.Cells.Item(88, 8).Value = 950.6667  # H88: was 996.8
$ws.Cells.Item(88, 10).Value = 995.8  # J88: was 1064.75
$ws.Cells.Item(88, 12).Value = 995.8  # L88: was 1064.75
$ws.Cells.Item(88, 14).Value = -1807.8  # N88: was -1876.75

$ws.Cells.Item(91, 8).Value = 950.6667  # H91: was 996.8
$ws.Cells.Item(91, 10).Value = 995.8  # J91: was 1064.75
$ws.Cells.Item(91, 12).Value = 995.8  # L91: was 1064.75
$ws.Cells.Item(91, 14).Value = -3803.8  # N91: was -3872.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 15903.333  # H2: was 6597.75
$ws.Cells.Item(2, 9).Value = 1105.5  # I2: was 1040.4286
$ws.Cells.Item(2, 11).Value = 1105.5  # K2: was 1040.4286
$ws.Cells.Item(2, 13).Value = -992.5  # M2: was -927.4286

$ws.Cells.Item(5, 8).Value = 129.4  # H5: was 162.33333
$ws.Cells.Item(5, 9).Value = 155.66667  # I5: was 196
$ws.Cells.Item(5, 10).Value = 90  # J5: was 95
$ws.Cells.Item(5, 11).Value = 155.66667  # K5: was 196
$ws.Cells.Item(5, 12).Value = 90  # L5: was 95
$ws.Cells.Item(5, 13).Value = -43.66667000000001  # M5: was -84
$ws.Cells.Item(5, 14).Value = -314  # N5: was -319

$ws.Cells.Item(45, 8).Value = 3327.4  # H45: was 3600.8462
$ws.Cells.Item(45, 9).Value = 2395  # I45: was 2733
$ws.Cells.Item(45, 11).Value = 2395  # K45: was 2733
$ws.Cells.Item(45, 13).Value = -2018  # M45: was -2356

$ws.Cells.Item(74, 8).Value = 3845.6  # H74: was 3729.9375
$ws.Cells.Item(74, 9).Value = 3798  # I74: was 3617.7
$ws.Cells.Item(74, 11).Value = 3798  # K74: was 3617.7
$ws.Cells.Item(74, 13).Value = -2924  # M74: was -2743.7

$ws.Cells.Item(77, 8).Value = 3845.6  # H77: was 3729.9375
$ws.Cells.Item(77, 9).Value = 3798  # I77: was 3617.7
$ws.Cells.Item(77, 11).Value = 18990  # K77: was 18088.5
$ws.Cells.Item(77, 13).Value = -14622  # M77: was -13720.5

$ws.Cells.Item(116, 8).Value = 15903.333  # H116: was 6597.75
$ws.Cells.Item(116, 9).Value = 1105.5  # I116: was 1040.4286
$ws.Cells.Item(116, 11).Value = 1105.5  # K116: was 1040.4286
$ws.Cells.Item(116, 13).Value = 1188.5  # M116: was 1253.5714

$ws.Cells.Item(122, 8).Value = 2186  # H122: was 2872
$ws.Cells.Item(122, 9).Value = 2186  # I122: was 2872
$ws.Cells.Item(122, 11).Value = 6558  # K122: was 8616
$ws.Cells.Item(122, 13).Value = -4108  # M122: was -6166

$ws.Cells.Item(132, 8).Value = 1212  # H132: was 1211.1428
$ws.Cells.Item(132, 9).Value = 1175.7  # I132: was 1177.909
$ws.Cells.Item(132, 11).Value = 3527.1  # K132: was 3533.727
$ws.Cells.Item(132, 13).Value = -997.1000000000004  # M132: was -1003.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 15903.333  # H3: was 6597.75
$ws.Cells.Item(3, 9).Value = 1105.5  # I3: was 1040.4286
$ws.Cells.Item(3, 11).Value = 1105.5  # K3: was 1040.4286
$ws.Cells.Item(3, 13).Value = -991.5  # M3: was -926.4286

$ws.Cells.Item(4, 8).Value = 129.4  # H4: was 162.33333
$ws.Cells.Item(4, 9).Value = 155.66667  # I4: was 196
$ws.Cells.Item(4, 10).Value = 90  # J4: was 95
$ws.Cells.Item(4, 11).Value = 155.66667  # K4: was 196
$ws.Cells.Item(4, 12).Value = 90  # L4: was 95
$ws.Cells.Item(4, 13).Value = -40.66667000000001  # M4: was -81
$ws.Cells.Item(4, 14).Value = -320  # N4: was -325

$ws.Cells.Item(20, 8).Value = 1545.2222  # H20: was 1733
$ws.Cells.Item(20, 9).Value = 532.75  # I20: was 669.3333
$ws.Cells.Item(20, 10).Value = 2355.2  # J20: was 2371.2
$ws.Cells.Item(20, 11).Value = 532.75  # K20: was 669.3333
$ws.Cells.Item(20, 12).Value = 2355.2  # L20: was 2371.2
$ws.Cells.Item(20, 13).Value = -285.75  # M20: was -422.3333
$ws.Cells.Item(20, 14).Value = -2849.2  # N20: was -2865.2

$ws.Cells.Item(80, 8).Value = 553.4666999999999  # H80: was 583.7857
$ws.Cells.Item(80, 9).Value = 280.22223  # I80: was 280.77777
$ws.Cells.Item(80, 10).Value = 963.3333  # J80: was 1129.2
$ws.Cells.Item(80, 11).Value = 280.22223  # K80: was 280.77777
$ws.Cells.Item(80, 12).Value = 963.3333  # L80: was 1129.2
$ws.Cells.Item(80, 13).Value = 717.7777699999999  # M80: was 717.2222300000001
$ws.Cells.Item(80, 14).Value = -2959.3333  # N80: was -3125.2

$ws.Cells.Item(83, 8).Value = 553.4666999999999  # H83: was 583.7857
$ws.Cells.Item(83, 9).Value = 280.22223  # I83: was 280.77777
$ws.Cells.Item(83, 10).Value = 963.3333  # J83: was 1129.2
$ws.Cells.Item(83, 11).Value = 1401.11115  # K83: was 1403.88885
$ws.Cells.Item(83, 12).Value = 4816.6665  # L83: was 5646
$ws.Cells.Item(83, 13).Value = 3590.88885  # M83: was 3588.11115
$ws.Cells.Item(83, 14).Value = -14800.6665  # N83: was -15630

$ws.Cells.Item(94, 8).Value = 0  # H94: was 398
$ws.Cells.Item(94, 9).Value = 0  # I94: was 475.2
$ws.Cells.Item(94, 10).Value = 0  # J94: was 205
$ws.Cells.Item(94, 11).Value = 0  # K94: was 475.2
$ws.Cells.Item(94, 12).ClearContents()  # L94: was 205
$ws.Cells.Item(94, 13).ClearContents()  # M94: was -24.19999999999999
$ws.Cells.Item(94, 14).Value = 0  # N94: was -1107

$ws.Cells.Item(105, 8).Value = 18183316  # H105: was 30304366
$ws.Cells.Item(105, 9).Value = 18183316  # I105: was 30304366
$ws.Cells.Item(105, 11).Value = 18183316  # K105: was 30304366
$ws.Cells.Item(105, 13).Value = -18181569  # M105: was -30302619

$ws.Cells.Item(112, 8).Value = 0  # H112: was 30000
$ws.Cells.Item(112, 10).Value = 0  # J112: was 30000
$ws.Cells.Item(112, 12).ClearContents()  # L112: was 30000
$ws.Cells.Item(112, 14).Value = 0  # N112: was -32954

$ws.Cells.Item(134, 8).Value = 1420.625  # H134: was 1494.5652
$ws.Cells.Item(134, 9).Value = 1111.8  # I134: was 1185.0526
$ws.Cells.Item(134, 11).Value = 3335.4  # K134: was 3555.1578
$ws.Cells.Item(134, 13).Value = -800.3999999999996  # M134: was -1020.1578

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(42, 8).Value = 20000  # H42: was 9000
$ws.Cells.Item(42, 9).Value = 20000  # I42: was 0
$ws.Cells.Item(42, 10).Value = 0  # J42: was 9000
$ws.Cells.Item(42, 11).Value = 20000  # K42: was 0
$ws.Cells.Item(42, 12).ClearContents()  # L42: was 9000
$ws.Cells.Item(42, 13).Value = -19407  # M42: was None
$ws.Cells.Item(42, 14).Value = 0  # N42: was -10186

$ws.Cells.Item(117, 8).Value = 42069  # H117: was 0
$ws.Cells.Item(117, 10).Value = 42069  # J117: was 0
$ws.Cells.Item(117, 12).Value = 42069  # L117: was 0
$ws.Cells.Item(117, 14).Value = -51247  # N117: was None

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 852  # H14: was 1401.3334
$ws.Cells.Item(14, 9).Value = 852  # I14: was 1401.3334
$ws.Cells.Item(14, 11).Value = 2556  # K14: was 4204.0002
$ws.Cells.Item(14, 13).Value = -2383  # M14: was -4031.0002

$ws.Cells.Item(55, 8).Value = 8073.9  # H55: was 7549
$ws.Cells.Item(55, 9).Value = 704.5  # I55: was 569.6667
$ws.Cells.Item(55, 10).Value = 9916.25  # J55: was 10166.25
$ws.Cells.Item(55, 11).Value = 2113.5  # K55: was 1709.0001
$ws.Cells.Item(55, 12).Value = 29748.75  # L55: was 30498.75
$ws.Cells.Item(55, 13).Value = -1936.5  # M55: was -1532.0001
$ws.Cells.Item(55, 14).Value = -30102.75  # N55: was -30852.75

$ws.Cells.Item(137, 8).Value = 2722  # H137: was 3184.8333
$ws.Cells.Item(137, 9).Value = 686.6667  # I137: was 865
$ws.Cells.Item(137, 10).Value = 5775  # J137: was 4344.75
$ws.Cells.Item(137, 11).Value = 2060.0001  # K137: was 2595
$ws.Cells.Item(137, 12).Value = 17325  # L137: was 13034.25
$ws.Cells.Item(137, 13).Value = 3039.9999  # M137: was 2505
$ws.Cells.Item(137, 14).Value = -27525  # N137: was -23234.25

$ws.Cells.Item(138, 8).Value = 6682.143  # H138: was 6733.3335
$ws.Cells.Item(138, 9).Value = 1800  # I138: was 1750
$ws.Cells.Item(138, 10).Value = 9394.444  # J138: was 9225
$ws.Cells.Item(138, 11).Value = 5400  # K138: was 5250
$ws.Cells.Item(138, 12).Value = 28183.332  # L138: was 27675
$ws.Cells.Item(138, 13).Value = -260  # M138: was -110
$ws.Cells.Item(138, 14).Value = -38463.33199999999  # N138: was -37955

$ws.Cells.Item(139, 8).Value = 4631.1665  # H139: was 5473.8
$ws.Cells.Item(139, 9).Value = 4631.1665  # I139: was 5473.8
$ws.Cells.Item(139, 11).Value = 13893.4995  # K139: was 16421.4
$ws.Cells.Item(139, 13).Value = -8753.499500000002  # M139: was -11281.4

$ws.Cells.Item(140, 8).Value = 3104.8333  # H140: was 3384.3635
$ws.Cells.Item(140, 9).Value = 2513.889  # I140: was 2824.375
$ws.Cells.Item(140, 11).Value = 7541.667  # K140: was 8473.125
$ws.Cells.Item(140, 13).Value = -2361.667  # M140: was -3293.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 45557.285  # H24: was 37777.4
$ws.Cells.Item(24, 10).Value = 45557.285  # J24: was 37777.4
$ws.Cells.Item(24, 12).Value = 45557.285  # L24: was 37777.4
$ws.Cells.Item(24, 14).Value = -45903.285  # N24: was -38123.4

$ws.Cells.Item(98, 8).Value = 9321.799999999999  # H98: was 9152.25
$ws.Cells.Item(98, 10).Value = 9321.799999999999  # J98: was 9152.25
$ws.Cells.Item(98, 12).Value = 9321.799999999999  # L98: was 9152.25
$ws.Cells.Item(98, 14).Value = -15311.8  # N98: was -15142.25

$ws.Cells.Item(102, 8).Value = 3569.7334  # H102: was 3762.4285
$ws.Cells.Item(102, 9).Value = 3116.4443  # I102: was 3397
$ws.Cells.Item(102, 11).Value = 3116.4443  # K102: was 3397
$ws.Cells.Item(102, 13).Value = -1494.4443  # M102: was -1775

$ws.Cells.Item(126, 8).Value = 2499.5  # H126: was 2498.5
$ws.Cells.Item(126, 9).Value = 2499  # I126: was 2498.5
$ws.Cells.Item(126, 10).Value = 2500  # J126: was 0
$ws.Cells.Item(126, 11).Value = 7497  # K126: was 7495.5
$ws.Cells.Item(126, 12).Value = 7500  # L126: was 0
$ws.Cells.Item(126, 13).Value = -5027  # M126: was -5025.5
$ws.Cells.Item(126, 14).Value = -12440  # N126: was None

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 7438.45  # H40: was 7519.222
$ws.Cells.Item(40, 9).Value = 6898.353  # I40: was 7146.6875
$ws.Cells.Item(40, 10).Value = 10499  # J40: was 10499.5
$ws.Cells.Item(40, 11).Value = 6898.353  # K40: was 7146.6875
$ws.Cells.Item(40, 12).Value = 10499  # L40: was 10499.5
$ws.Cells.Item(40, 13).Value = -6762.353  # M40: was -7010.6875
$ws.Cells.Item(40, 14).Value = -10771  # N40: was -10771.5

$ws.Cells.Item(122, 8).Value = 4856.143  # H122: was 5710.625
$ws.Cells.Item(122, 9).Value = 5198.8335  # I122: was 6126.4287
$ws.Cells.Item(122, 11).Value = 15596.5005  # K122: was 18379.2861
$ws.Cells.Item(122, 13).Value = -13146.5005  # M122: was -15929.2861

$ws.Cells.Item(136, 8).Value = 2998.25  # H136: was 2999
$ws.Cells.Item(136, 9).Value = 2998.25  # I136: was 2999
$ws.Cells.Item(136, 11).Value = 8994.75  # K136: was 8997
$ws.Cells.Item(136, 13).Value = -6444.75  # M136: was -6447

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(31, 8).Value = 15000  # H31: was 0
$ws.Cells.Item(31, 10).Value = 15000  # J31: was 0
$ws.Cells.Item(31, 12).Value = 15000  # L31: was 0
$ws.Cells.Item(31, 14).Value = -15696  # N31: was None

$ws.Cells.Item(100, 8).Value = 773.25  # H100: was 597.3333
$ws.Cells.Item(100, 9).Value = 864.3333  # I100: was 616.8
$ws.Cells.Item(100, 11).Value = 1728.6666  # K100: was 1233.6
$ws.Cells.Item(100, 13).Value = -1187.6666  # M100: was -692.5999999999999

$ws.Cells.Item(107, 8).Value = 55556544  # H107: was 22223046
$ws.Cells.Item(107, 9).Value = 111111650  # I107: was 27778444
$ws.Cells.Item(107, 10).Value = 1446.6666  # J107: was 1450
$ws.Cells.Item(107, 11).Value = 333334950  # K107: was 83335332
$ws.Cells.Item(107, 12).Value = 4339.9998  # L107: was 4350
$ws.Cells.Item(107, 13).Value = -333333030  # M107: was -83333412
$ws.Cells.Item(107, 14).Value = -8179.9998  # N107: was -8190

$ws.Cells.Item(113, 8).Value = 754.63635  # H113: was 945.1
$ws.Cells.Item(113, 9).Value = 460  # I113: was 750
$ws.Cells.Item(113, 10).Value = 1000.1667  # J113: was 1140.2
$ws.Cells.Item(113, 11).Value = 1380  # K113: was 2250
$ws.Cells.Item(113, 12).Value = 3000.5001  # L113: was 3420.6
$ws.Cells.Item(113, 13).Value = 790  # M113: was -80
$ws.Cells.Item(113, 14).Value = -7340.5001  # N113: was -7760.6

$ws.Cells.Item(136, 8).Value = 2976.5625  # H136: was 3076.1333
$ws.Cells.Item(136, 9).Value = 2077.9167  # I136: was 2132
$ws.Cells.Item(136, 11).Value = 6233.750100000001  # K136: was 6396
$ws.Cells.Item(136, 13).Value = -3683.750100000001  # M136: was -3846
